$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data scraped on Fri Nov 29 22:56:57 UTC 2024

$ws.Range("D2").Value = "97.420.13"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "3.600.36"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "244.51"
$ws.Range("E5").Value = "  +2.77%  "
$ws.Range("E6").Value = "  +17.18%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("D8").Value = "0.427"
$ws.Range("E8").Value = "  +5.84%  "
$ws.Range("E9").Value = "  -0.07%  "
$ws.Range("E10").Value = "  +1.94%  "
$ws.Range("D11").Value = "3.598.66"
$ws.Range("E11").Value = "  +0.26%  "
$ws.Range("D12").Value = "44.74"
$ws.Range("E12").Value = "  +3.85%  "
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -0.01%  "
$ws.Range("D15").Value = "4.266.37"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "97.352.23"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("E17").Value = "  +2.51%  "
$ws.Range("D18").Value = "3.599.91"
$ws.Range("E18").Value = "  +0.31%  "
$ws.Range("E19").Value = "  +0.31%  "
$ws.Range("D20").Value = "12.68"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("D21").Value = "18.25"
$ws.Range("E21").Value = "  +1.19%  "
$ws.Range("D22").Value = "0.527"
$ws.Range("E22").Value = "  +7.69%  "
$ws.Range("B23").Value = "SuiNetwork"
$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").Value = "3.49"
$ws.Range("E23").Value = "  +0.29%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D24").Value = "519.11"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.00"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").Value = "103.54"
$ws.Range("E27").Value = "  +7.88%  "
$ws.Range("E28").Value = "  +3.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.180"
$ws.Range("E29").Value = "  +22.78%  "
$ws.Range("D30").Value = "3.794.34"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("E31").Value = "  -2.48%  "
$ws.Range("D32").Value = "11.97"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  +6.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  -0.02%  "
$ws.Range("D36").Value = "31.87"
$ws.Range("E36").Value = "  -0.89%  "
$ws.Range("E37").Value = "  +3.60%  "
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").Value = "618.67"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").Value = "8.78"
$ws.Range("E40").Value = "  +0.47%  "
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "0.465"
$ws.Range("E45").Value = "  +46.46%  "
$ws.Range("B46").Value = "Filecoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D46").Value = "6.12"
$ws.Range("E46").Value = "  +5.75%  "
$ws.Range("E47").Value = "  +6.54%  "
$ws.Range("E48").Value = "  +1.40%  "
$ws.Range("E49").Value = "  +0.99%  "
$ws.Range("D50").Value = "8.68"
$ws.Range("E50").Value = "  +5.41%  "
$ws.Range("D51").Value = "33.11"
$ws.Range("E51").Value = "  -5.81%  "
